$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.113.93'
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").Value = '1.787.56'
$ws.Range("E3").Value = '  -2.99%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''224.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("E6").Value = '  -1.50%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '''32.90'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.52%  '
$ws.Range("E9").Value = '  -2.91%  '
$ws.Range("E10").Value = '  -1.15%  '
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").Value = '2.044.29'
$ws.Range("E12").Value = '  -3.23%  '
$ws.Range("D13").Value = '1.784.27'
$ws.Range("E13").Value = '  -3.22%  '
$ws.Range("D14").Value = '''10.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("E15").Value = '  -4.02%  '
$ws.Range("D16").Value = '34.048.99'
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("E17").Value = '  -5.27%  '
$ws.Range("D18").Value = '''67.89'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.83%  '
$ws.Range("D19").Value = '''245.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.62%  '
$ws.Range("D20").Value = '0.0₃0789'
$ws.Range("E20").Value = '  -1.85%  '
$ws.Range("D21").Value = '''0.998'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = '''10.81'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.51%  '
$ws.Range("E23").Value = '  -4.62%  '
$ws.Range("E24").Value = '  -2.90%  '
$ws.Range("D25").Value = '''160.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.79%  '
$ws.Range("D26").Value = '''16.35'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.47%  '
$ws.Range("D27").Value = '''7.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.94%  '
$ws.Range("E28").Value = '  -2.68%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("E30").Value = '  -4.70%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("E32").Value = '  -4.29%  '
$ws.Range("D33").Value = '''3.51'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.27%  '
$ws.Range("E34").Value = '  -6.49%  '
$ws.Range("D35").Value = '1.397.29'
$ws.Range("E35").Value = '  -4.63%  '
$ws.Range("E36").Value = '  -1.22%  '
$ws.Range("E37").Value = '  -1.96%  '
$ws.Range("E38").Value = '  -4.00%  '
$ws.Range("E39").Value = '  +2.93%  '
$ws.Range("E40").Value = '  -0.47%  '
$ws.Range("D41").Value = '''0.915'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.92%  '
$ws.Range("D42").Value = '''2.69'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.09%  '
$ws.Range("D43").Value = '''78.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.37%  '
$ws.Range("E44").Value = '  +14.36%  '
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").Value = '''0.0498'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.28%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '''12.45'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.51%  '
$ws.Range("D48").Value = '''107.86'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.72%  '
$ws.Range("E49").Value = '  -4.61%  '
$ws.Range("D50").Value = '1.943.69'
$ws.Range("E50").Value = '  -3.15%  '
$ws.Range("E51").Value = '  -0.47%  '
